# update figures 1 and 2
#
# The "OLD" sheet (flowchart source data) gets a new row inserted right
# before the existing "Placebo/sham use" row under
# Study-related factors > Study treatment, for a new subtheme
# "Medication administration" (frequency left blank).  Every row that used
# to follow (old rows 25-31) shifts down by one (new rows 26-32).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OLD")

# Insert a new row 25, pushing the old row 25 ("Placebo/sham use") and
# everything below it down by one row.
$ws.Rows(25).Insert()

$ws.Range("A25").Value = "Study-related factors"
$ws.Range("B25").Value = "Study treatment"
$ws.Range("C25").Value = "Medication administration"

# Leave the final click state matching the source workbook.
$ws.Range("F14").Select()
